# Generate Report for Handback
# Adds a new handed-back file (6920636f-3f4a-43ea-93de-50b787ffc917.md) as a
# new row to the Overview, zh-cn and de-de report tables.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# Expand each table by one row - this keeps the table ref / autoFilter ref
# and worksheet dimension in sync automatically.
$lo1 = $ws1.ListObjects.Item(1)
$lo2 = $ws2.ListObjects.Item(1)
$lo3 = $ws3.ListObjects.Item(1)

$lo1.ListRows.Add() | Out-Null
$lo2.ListRows.Add() | Out-Null
$lo3.ListRows.Add() | Out-Null

$fileName   = "6920636f-3f4a-43ea-93de-50b787ffc917.md"
$pathName   = "e2e\6920636f-3f4a-43ea-93de-50b787ffc917.md"
$extension  = ".md"
$statusMsg  = "Handed back: in sync with en-US"

$zhXlf = "6920636f-3f4a-43ea-93de-50b787ffc917.1b10f6bddb18522f912d09ef8f508f4ad125f188.zh-cn.xlf"
$deXlf = "6920636f-3f4a-43ea-93de-50b787ffc917.1b10f6bddb18522f912d09ef8f508f4ad125f188.de-de.xlf"

$zhHandoffDate  = "2017-02-21 03:50:24"
$zhHandbackDate = "2017-02-21 03:51:18"
$deHandoffDate  = "2017-02-21 03:50:39"
$deHandbackDate = "2017-02-21 03:51:41"

$overviewDate = "2017-02-21 03:50:39"

$sourceUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/a2994b5124b31d56cbc9145f18983b9ceea72dad/e2e/"
$zhUrlBase     = "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/9a7b49d1060f72a912ea4de931e4e60ee77a01c3/e2e/"
$deUrlBase     = "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/c48f5a270b39bb019d54ad7ced67646d88150e2f/e2e/"

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (columns A..G, new row 3)
# ---------------------------------------------------------------------------
$ws1.Range("A3").Value = $fileName
$ws1.Range("C3").Value = $extension
$ws1.Range("E3").Value = $statusMsg
$ws1.Range("F3").Value = $statusMsg
$ws1.Range("G3").Value = $overviewDate
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B3"), ($sourceUrlBase + $fileName), [Type]::Missing, [Type]::Missing, $pathName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn  (columns A..R, new row 3)
# ---------------------------------------------------------------------------
$ws2.Range("B3").Value = $extension
$ws2.Range("C3").Value = $statusMsg
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'True"
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("H3").Value = $zhHandoffDate
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K3").Value = $zhXlf
$ws2.Range("L3").Value = $zhHandbackDate
$ws2.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("O3").Value = "'True"
$ws2.Range("Q3").Value = "'False"

$ws2.Hyperlinks.Add($ws2.Range("A3"), ($sourceUrlBase + $fileName), [Type]::Missing, [Type]::Missing, $fileName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("J3"), ($zhUrlBase + $fileName), [Type]::Missing, [Type]::Missing, $fileName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de  (columns A..R, new row 3)
# ---------------------------------------------------------------------------
$ws3.Range("B3").Value = $extension
$ws3.Range("C3").Value = $statusMsg
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'True"
$ws3.Range("G3").Value = $deXlf
$ws3.Range("H3").Value = $deHandoffDate
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K3").Value = $deXlf
$ws3.Range("L3").Value = $deHandbackDate
$ws3.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("O3").Value = "'True"
$ws3.Range("Q3").Value = "'False"

$ws3.Hyperlinks.Add($ws3.Range("A3"), ($sourceUrlBase + $fileName), [Type]::Missing, [Type]::Missing, $fileName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("J3"), ($deUrlBase + $fileName), [Type]::Missing, [Type]::Missing, $fileName) | Out-Null
